$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16: student #13 — homework 5 (G) grade entered ---
$ws.Range("G16").Value = 5

# --- Row 20: student #17 — homework 4 (F) grade corrected 0 -> 5,
#     and the stale "missing homework" green highlight is cleared ---
$ws.Range("H20").Copy() | Out-Null
$ws.Range("F20").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("F20").Value = 5

# --- Row 21: student #18 — homework 2 (D) grade corrected 0 -> 5,
#     and the stale "missing homework" green highlight is cleared ---
$ws.Range("H21").Copy() | Out-Null
$ws.Range("D21").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("D21").Value = 5

$excel.CutCopyMode = 0

# --- Update the on-screen view: scroll the frozen pane up and move the
#     active selection in the bottom-right pane to G15 ---
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("G15").Select() | Out-Null
